$wb = $excel.ActiveWorkbook

# --- Update "Publisher Max Sales" sheet (sheet1) ---
$wsMax = $wb.Worksheets.Item("Publisher Max Sales")

# Row 5 (Net sales): Q5,R5,S5 75000 -> 100000
$wsMax.Range("Q5").Value = 100000
$wsMax.Range("R5").Value = 100000
$wsMax.Range("S5").Value = 100000

# Row 21 (Investment Recoup): Q21,R21 75000 -> 100000, S21 cleared (was 50000)
$wsMax.Range("Q21").Value = 100000
$wsMax.Range("R21").Value = 100000
$wsMax.Range("S21").ClearContents()

# Row 22 (Revenue Split): S22 12500 -> 50000
$wsMax.Range("S22").Value = 50000

# Make "Publisher Max Sales" the active sheet/tab, with T21 selected
$wsMax.Activate()
$wsMax.Range("T21").Select()

$wb.Save()
